# Update the "version" tracking sheet with the latest component versions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- libexpat: R_2_2_5 -> R_2_2_6 (row 11) ---
$ws.Range("C11").Value = "R_2_2_6"
$ws.Range("D11").Value = "0816 R_2-2-5 => R_2_2_6"

# --- openssl: OpenSSL_1_1_0h -> OpenSSL_1_1_0i (row 23) ---
$ws.Range("B23").Value = "msvc15 / msvc15-OpenSSL_1_1_0i"
$ws.Range("C23").Value = "OpenSSL_1_1_0i"
$ws.Range("D23").Value = "0327 OpenSSL_1_1_0g => OpenSSL_1_1_0h`n0816 OpenSSL_1_1_0h => OpenSSL_1_1_0i"

# --- php: 7.2.8 -> 7.2.9 (row 28) ---
$ws.Range("C28").Value = "7.2.9"
$ws.Range("D28").Value = " * 0301`n0327 7.2.4`n0425 7.2.5`n0612 7.2.6`n0630 7.2.7`n0724 7.2.8`n0816 7.2.9"

# --- memcached: 1.5.9 -> 1.5.10 (row 37) ---
$ws.Range("B37").Value = "msvc15 / msvc15-1.5.10"
$ws.Range("C37").Value = "1.5.10"
$ws.Range("D37").Value = "0612 : 1.5.7 => 1.5.8`n0724 : 1.5.8 => 1.5.9`n0816 : 1.5.9 => 1.5.10 "

# Re-fit the rows whose wrapped text grew by one line.
$ws.Rows.AutoFit()

# --- Update the view: scroll position, zoom, and active selection ---
$win = $excel.ActiveWindow
$win.Zoom = 70
$ws.Range("A16").Select()
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("E22").Select()

$wb.Save()
